$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with new ROI data

$ws.Range("A2").Value = "EvapDOdt1"
$ws.Range("B2").Value = 334
$ws.Range("C2").Value = 1668
$ws.Range("D2").Value = 530
$ws.Range("E2").Value = 1970
$ws.Range("F2").Value = 2160
$ws.Range("G2").Value = 2560
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = "[]"
$ws.Range("J2").Value = "[1 1]"
$ws.Range("K2").Value = "[100 100]"

$ws.Range("A5").Value = "Bec"
$ws.Range("B5").Value = 911
$ws.Range("C5").Value = 1183
$ws.Range("D5").Value = 1235
$ws.Range("E5").Value = 1511
$ws.Range("F5").Value = 2160
$ws.Range("G5").Value = 2560
$ws.Range("H5").Value = 2.2999999999999998
$ws.Range("I5").Value = "[]"
$ws.Range("J5").Value = "[1 1]"
$ws.Range("K5").Value = "[100 100]"

$ws.Range("A7").Value = "NiLattice"
$ws.Range("B7").Value = 778
$ws.Range("C7").Value = 1372
$ws.Range("D7").Value = 1286
$ws.Range("E7").Value = 1442
$ws.Range("F7").Value = 2160
$ws.Range("G7").Value = 2560
$ws.Range("H7").Value = 2.2999999999999998
$ws.Range("I7").Value = "[]"
$ws.Range("J7").Value = "[1 1]"
$ws.Range("K7").Value = "[100 100]"

$ws.Range("A8").Value = "NILatticeKd"
$ws.Range("B8").Value = 708
$ws.Range("C8").Value = 1400
$ws.Range("D8").Value = 1346
$ws.Range("E8").Value = 1400
$ws.Range("F8").Value = 2160
$ws.Range("G8").Value = 2560
$ws.Range("H8").Value = 2.3000000000000003
$ws.Range("I8").Value = "[]"
$ws.Range("J8").Value = "[1 1]"
$ws.Range("K8").Value = "[100 100]"

$ws.Range("A14").Value = "NiLatticeDepthCalib"
$ws.Range("B14").Value = 759
$ws.Range("C14").Value = 1313
$ws.Range("D14").Value = 1236
$ws.Range("E14").Value = 1556
$ws.Range("F14").Value = 2160
$ws.Range("G14").Value = 2560
$ws.Range("H14").Value = 2.2999999999999998
$ws.Range("I14").Value = "[998 1357 130 100]"
$ws.Range("J14").Value = "[3 1]"
$ws.Range("K14").Value = "[170 100]"

$ws.Range("A15").Value = "BMPDloopTof3000"
$ws.Range("B15").Value = 490
$ws.Range("C15").Value = 1352
$ws.Range("D15").Value = 1341
$ws.Range("E15").Value = 1483
$ws.Range("F15").Value = 2160
$ws.Range("G15").Value = 2560
$ws.Range("H15").Value = 2.2999999999999998
$ws.Range("I15").Value = "[857 1387 260 75]"
$ws.Range("J15").Value = "[2 1]"
$ws.Range("K15").Value = "[560 100]"

$ws.Range("A24").Value = "NiLatticeBoBm"
$ws.Range("B24").Value = 510
$ws.Range("C24").Value = 1468
$ws.Range("D24").Value = 1333
$ws.Range("E24").Value = 1459
$ws.Range("F24").Value = 2160
$ws.Range("G24").Value = 2560
$ws.Range("H24").Value = 2.2999999999999998
$ws.Range("I24").Value = "[]"
$ws.Range("J24").Value = "[1 1]"
$ws.Range("K24").Value = "[100 100]"

$ws.Range("A25").Value = "BMPDloopTof5000"
$ws.Range("B25").Value = 454
$ws.Range("C25").Value = 1500
$ws.Range("D25").Value = 1365
$ws.Range("E25").Value = 1523
$ws.Range("F25").Value = 2160
$ws.Range("G25").Value = 2560
$ws.Range("H25").Value = 2.2999999999999998
$ws.Range("I25").Value = "[941 1405 180 100]"
$ws.Range("J25").Value = "[2 1]"
$ws.Range("K25").Value = "[750 100]"

$ws.Range("A26").Value = "NiLatticeDepthCalibTof6000"
$ws.Range("B26").Value = 634
$ws.Range("C26").Value = 1422
$ws.Range("D26").Value = 1333
$ws.Range("E26").Value = 1563
$ws.Range("F26").Value = 2160
$ws.Range("G26").Value = 2560
$ws.Range("H26").Value = 2.2999999999999998
$ws.Range("I26").Value = "[967 1410 130 100]"
$ws.Range("J26").Value = "[3 1]"
$ws.Range("K26").Value = "[270 100]"

$ws.Range("A27").Value = "NiLatticeDepthCalibTof5000"
$ws.Range("B27").Value = 709
$ws.Range("C27").Value = 1419
$ws.Range("D27").Value = 1287
$ws.Range("E27").Value = 1445
$ws.Range("F27").Value = 2160
$ws.Range("G27").Value = 2560
$ws.Range("H27").Value = 2.2999999999999998
$ws.Range("I27").Value = "[1024 1328 130 120]"
$ws.Range("J27").Value = "[3 1]"
$ws.Range("K27").Value = "[270 100]"

$ws.Range("A28").Value = "NiLatticeBoBmTof5000"
$ws.Range("B28").Value = 582
$ws.Range("C28").Value = 1484
$ws.Range("D28").Value = 1385
$ws.Range("E28").Value = 1505
$ws.Range("F28").Value = 2160
$ws.Range("G28").Value = 2560
$ws.Range("H28").Value = 2.2999999999999998
$ws.Range("I28").Value = "[]"
$ws.Range("J28").Value = "[1 1]"
$ws.Range("K28").Value = "[100 100]"

$ws.Range("A29").Value = "HfBecFullTof"
$ws.Range("B29").Value = 736
$ws.Range("C29").Value = 1358
$ws.Range("D29").Value = 1060
$ws.Range("E29").Value = 1686
$ws.Range("F29").Value = 2160
$ws.Range("G29").Value = 2560
$ws.Range("H29").Value = 2.2999999999999998
$ws.Range("I29").Value = "[]"
$ws.Range("J29").Value = "[1 1]"
$ws.Range("K29").Value = "[100 100]"
